$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Data")

# Row 4: bump row height to match the "thick bottom" look already used by the rows below it
$ws.Range("A4:L4").RowHeight = 15.6

# Row 5: replace the content with the new bulk-order test row, matching the style of rows 6-12.
# Copy row 6's formatting onto row 5 first (keeps the thick-bottom row style + per-column styles),
# then overwrite the cell values with the new test data.
$ws.Range("A6:L6").Copy()
$ws.Range("A5:L5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C5").Value = "MA70"
$ws.Range("D5").Value = "CA"
$ws.Range("E5").Value = 10066860
$ws.Range("F5").Value = "JA 0911"
$ws.Range("G5").Value = "OT"
$ws.Range("H5").Value = "S"
$ws.Range("I5").Value = "Header Comment 1"
$ws.Range("J5").Value = "Header 1"
$ws.Range("K5").Value = "4353CD"
$ws.Range("L5").Value = 1

# Update the active selection to H5 (matches the workbook's saved cursor position)
$ws.Range("H5").Select()
